# repull data, push all data, mean calculation
# Update column F (dSF) values on Sheet1 to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    4  = -6
    5  = -2
    6  = 0
    7  = -11
    8  = -2
    9  = 2
    10 = 2
    11 = -1
    12 = 3
    13 = 6
    14 = -5
    15 = -5
    16 = 3
    17 = 4
    18 = -1
    21 = 8
    22 = -3
    23 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
